$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data holds prices/percentages as literal text (t="inlineStr")
# rather than numbers. Stamping the Text number format ("@") before the
# assignment keeps Excel from auto-coercing these numeric-looking / percent-
# looking strings into actual numbers, so the cell keeps its original text type.
$cells = @(
    @{Addr='D2'; Val='309.37'}
    @{Addr='E2'; Val='-2.68%'}
    @{Addr='D3'; Val='37.27'}
    @{Addr='E3'; Val='-6.35%'}
    @{Addr='D4'; Val='5.114'}
    @{Addr='E4'; Val='-0.48%'}
    @{Addr='D5'; Val='0.07864'}
    @{Addr='E5'; Val='-4.07%'}
    @{Addr='D6'; Val='1.967'}
    @{Addr='E6'; Val='-5.17%'}
    @{Addr='D7'; Val='4.370'}
    @{Addr='E7'; Val='1.67%'}
    @{Addr='D8'; Val='8.284'}
    @{Addr='E8'; Val='-0.34%'}
    @{Addr='D9'; Val='3.101'}
    @{Addr='E9'; Val='-2.37%'}
    @{Addr='D10'; Val='0.9250'}
    @{Addr='E10'; Val='-1.00%'}
    @{Addr='D11'; Val='0.1338'}
    @{Addr='E11'; Val='-1.46%'}
    @{Addr='D12'; Val='0.1943'}
    @{Addr='E12'; Val='-1.94%'}
    @{Addr='D13'; Val='0.08955'}
    @{Addr='E13'; Val='-1.35%'}
    @{Addr='D14'; Val='0.03439'}
    @{Addr='E14'; Val='-1.27%'}
    @{Addr='E15'; Val='-1.10%'}
    @{Addr='D16'; Val='0.001385'}
    @{Addr='E16'; Val='-1.14%'}
    @{Addr='D17'; Val='0.005864'}
    @{Addr='E17'; Val='-8.11%'}
    @{Addr='D18'; Val='3.590'}
    @{Addr='E18'; Val='-2.52%'}
    @{Addr='D19'; Val='0.3400'}
    @{Addr='E19'; Val='-2.04%'}
    @{Addr='D20'; Val='0.1294'}
    @{Addr='E20'; Val='0.09%'}
    @{Addr='D21'; Val='5.010'}
    @{Addr='E21'; Val='2.14%'}
    @{Addr='D22'; Val='0.2489'}
    @{Addr='E22'; Val='1.56%'}
    @{Addr='D23'; Val='0.02105'}
    @{Addr='E23'; Val='5,162.71%'}
    @{Addr='D24'; Val='0.04346'}
    @{Addr='E24'; Val='0.43%'}
    @{Addr='E25'; Val='-0.68%'}
    @{Addr='D26'; Val='0.004538'}
    @{Addr='E26'; Val='-4.39%'}
    @{Addr='D27'; Val='0.0001352'}
    @{Addr='E27'; Val='4.01%'}
    @{Addr='D39'; Val='0.02294'}
    @{Addr='E39'; Val='3.39%'}
    @{Addr='D40'; Val='0.05023'}
    @{Addr='E40'; Val='-3.82%'}
    @{Addr='D41'; Val='0.007664'}
    @{Addr='E41'; Val='0.44%'}
    @{Addr='D42'; Val='0.009840'}
    @{Addr='E42'; Val='1.75%'}
    @{Addr='D43'; Val='0.1353'}
    @{Addr='E43'; Val='-2.16%'}
    @{Addr='D44'; Val='0.002064'}
    @{Addr='E44'; Val='-3.13%'}
    @{Addr='D45'; Val='0.008412'}
    @{Addr='E45'; Val='-8.57%'}
    @{Addr='D46'; Val='0.00006788'}
    @{Addr='E46'; Val='3.53%'}
    @{Addr='D47'; Val='0.00000000750'}
    @{Addr='E47'; Val='0.03%'}
    @{Addr='D48'; Val='0.003006'}
    @{Addr='E48'; Val='8.19%'}
    @{Addr='D49'; Val='0.001301'}
    @{Addr='E49'; Val='8.39%'}
    @{Addr='D50'; Val='0.00002101'}
    @{Addr='E50'; Val='0.03%'}
    @{Addr='D51'; Val='0.0002001'}
    @{Addr='E51'; Val='0.03%'}
)

foreach ($cell in $cells) {
    $rng = $ws.Range($cell.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.Val
}
